# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# The workbook's roll-up sheet "总计" becomes the new "2022-Q1" per-quarter
# holdings sheet (same sheetId/r:id it already had), and a brand new "总计"
# sheet is appended after it containing the refreshed roll-up (old rows
# shifted down one, with a new first data row for 2022-Q1).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# "2021-Q4" already carries the exact header / index-column formatting
# (bold, thin-bordered, centered) that every quarter sheet uses -- reuse it
# via copy/paste-format instead of fabricating a new style, so the shared
# style table isn't perturbed.
$styleDonor = $wb.Worksheets.Item("2021-Q4")

# Helper: write a value as genuine text (the source data stores numeric-
# looking figures like "0.64" as text, not numbers) without perturbing the
# cell's style.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

function Copy-HeaderFormat($destRange) {
    $styleDonor.Range("B1").Copy() | Out-Null
    $destRange.PasteSpecial(-4122) | Out-Null
}

function Copy-IndexFormat($destRange) {
    $styleDonor.Range("A2").Copy() | Out-Null
    $destRange.PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# 1. Rename the current "总计" sheet to "2022-Q1" -- it keeps its original
#    sheetId/r:id, matching the diff (sheetId="5" r:id="rId5").
# ---------------------------------------------------------------------------
$q1_2022 = $wb.Worksheets.Item("总计")
$q1_2022.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# 2. Replace its contents with the 2022-Q1 fund-holding breakdown.
# ---------------------------------------------------------------------------
$q1_2022.Cells.Clear()

Set-TextValue $q1_2022.Range("B1") "基金代码"
Set-TextValue $q1_2022.Range("C1") "基金名称"
Set-TextValue $q1_2022.Range("D1") "基金规模"
Set-TextValue $q1_2022.Range("E1") "股票总仓位"
Set-TextValue $q1_2022.Range("F1") "仓位占比"
Set-TextValue $q1_2022.Range("G1") "持有市值(亿元)"
Set-TextValue $q1_2022.Range("H1") "仓位排名"
Copy-HeaderFormat $q1_2022.Range("B1:H1")

$q1_2022.Range("A2").Value = 0
Set-TextValue $q1_2022.Range("B2") "970020"
Set-TextValue $q1_2022.Range("C2") "信达价值精选一年持有期灵活配置混合A"
Set-TextValue $q1_2022.Range("D2") "0.64"
Set-TextValue $q1_2022.Range("E2") "56.02"
Set-TextValue $q1_2022.Range("F2") "6.59"
Set-TextValue $q1_2022.Range("G2") "0.0422"
$q1_2022.Range("H2").Value = 1

$q1_2022.Range("A3").Value = 1
Set-TextValue $q1_2022.Range("B3") "970021"
Set-TextValue $q1_2022.Range("C3") "信达价值精选一年持有期灵活配置混合B"
Set-TextValue $q1_2022.Range("D3") "0.53"
Set-TextValue $q1_2022.Range("E3") "56.02"
Set-TextValue $q1_2022.Range("F3") "6.59"
Set-TextValue $q1_2022.Range("G3") "0.0349"
$q1_2022.Range("H3").Value = 1

Copy-IndexFormat $q1_2022.Range("A2:A3")

# ---------------------------------------------------------------------------
# 3. Append a brand new "总计" sheet after "2022-Q1" with the refreshed
#    roll-up (sheetId="6" r:id="rId6").
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1_2022)
$total.Name = "总计"

Set-TextValue $total.Range("B1") "日期"
Set-TextValue $total.Range("C1") "持有数量(只)"
Set-TextValue $total.Range("D1") "持有市值(亿元)"
Copy-HeaderFormat $total.Range("B1:D1")

$total.Range("A2").Value = 0
Set-TextValue $total.Range("B2") "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.08

$total.Range("A3").Value = 1
Set-TextValue $total.Range("B3") "2021-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.05

$total.Range("A4").Value = 2
Set-TextValue $total.Range("B4") "2021-Q3"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.03

$total.Range("A5").Value = 3
Set-TextValue $total.Range("B5") "2021-Q2"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 0.88

$total.Range("A6").Value = 4
Set-TextValue $total.Range("B6") "2021-Q1"
$total.Range("C6").Value = 5
$total.Range("D6").Value = 0.9

Copy-IndexFormat $total.Range("A2:A6")
